$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4145.161
$ws.Range("J76").Value = 4183.3335
$ws.Range("L76").Value = 4183.3335
$ws.Range("N76").Value = -4813.3335
$ws.Range("H79").Value = 4145.161
$ws.Range("J79").Value = 4183.3335
$ws.Range("L79").Value = 4183.3335
$ws.Range("N79").Value = -6367.3335
$ws.Range("H96").Value = 784.93335
$ws.Range("I96").Value = 418.44446
$ws.Range("J96").Value = 1334.6666
$ws.Range("K96").Value = 1255.33338
$ws.Range("L96").Value = 4003.9998
$ws.Range("M96").Value = 117.66662
$ws.Range("N96").Value = -6749.9998
$ws.Range("H141").Value = 1954.7
$ws.Range("I141").Value = 1506.7858
$ws.Range("J141").Value = 2999.8333
$ws.Range("K141").Value = 4520.357400000001
$ws.Range("L141").Value = 8999.499899999999
$ws.Range("M141").Value = 659.6425999999992
$ws.Range("N141").Value = -19359.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1593.9286
$ws.Range("I2").Value = 1416.8334
$ws.Range("J2").Value = 2656.5
$ws.Range("K2").Value = 1416.8334
$ws.Range("L2").Value = 2656.5
$ws.Range("M2").Value = -1303.8334
$ws.Range("N2").Value = -2882.5
$ws.Range("H45").Value = 1972.4231
$ws.Range("I45").Value = 1723.9231
$ws.Range("J45").Value = 2220.923
$ws.Range("K45").Value = 1723.9231
$ws.Range("L45").Value = 2220.923
$ws.Range("M45").Value = -1346.9231
$ws.Range("N45").Value = -2974.923
$ws.Range("H74").Value = 5190.625
$ws.Range("I74").Value = 5798.05
$ws.Range("J74").Value = 2153.5
$ws.Range("K74").Value = 5798.05
$ws.Range("L74").Value = 2153.5
$ws.Range("M74").Value = -4924.05
$ws.Range("N74").Value = -3901.5
$ws.Range("H77").Value = 5190.625
$ws.Range("I77").Value = 5798.05
$ws.Range("J77").Value = 2153.5
$ws.Range("K77").Value = 28990.25
$ws.Range("L77").Value = 10767.5
$ws.Range("M77").Value = -24622.25
$ws.Range("N77").Value = -19503.5
$ws.Range("H116").Value = 1593.9286
$ws.Range("I116").Value = 1416.8334
$ws.Range("J116").Value = 2656.5
$ws.Range("K116").Value = 1416.8334
$ws.Range("L116").Value = 2656.5
$ws.Range("M116").Value = 877.1666
$ws.Range("N116").Value = -7244.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1593.9286
$ws.Range("I3").Value = 1416.8334
$ws.Range("J3").Value = 2656.5
$ws.Range("K3").Value = 1416.8334
$ws.Range("L3").Value = 2656.5
$ws.Range("M3").Value = -1302.8334
$ws.Range("N3").Value = -2884.5
$ws.Range("H94").Value = 2381.6
$ws.Range("I94").Value = 2208
$ws.Range("J94").Value = 2425
$ws.Range("K94").Value = 2208
$ws.Range("L94").Value = 2425
$ws.Range("M94").Value = -1757
$ws.Range("N94").Value = -3327
$ws.Range("H134").Value = 2043.7949
$ws.Range("I134").Value = 1924.4138
$ws.Range("J134").Value = 2390
$ws.Range("K134").Value = 5773.2414
$ws.Range("L134").Value = 7170
$ws.Range("M134").Value = -3238.2414
$ws.Range("N134").Value = -12240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2178.24
$ws.Range("I132").Value = 1584.9412
$ws.Range("J132").Value = 3439
$ws.Range("K132").Value = 4754.8236
$ws.Range("L132").Value = 10317
$ws.Range("M132").Value = -2224.8236
$ws.Range("N132").Value = -15377

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224
$ws.Range("H131").Value = 762.0769
$ws.Range("I131").Value = 344.25
$ws.Range("J131").Value = 947.7778
$ws.Range("K131").Value = 1032.75
$ws.Range("L131").Value = 2843.3334
$ws.Range("M131").Value = 4007.25
$ws.Range("N131").Value = -12923.3334
$ws.Range("H139").Value = 4169263.2
$ws.Range("I139").Value = 8333933.5
$ws.Range("J139").Value = 4593.3335
$ws.Range("K139").Value = 25001800.5
$ws.Range("L139").Value = 13780.0005
$ws.Range("M139").Value = -24996660.5
$ws.Range("N139").Value = -24060.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5889.7856
$ws.Range("I70").Value = 6137.25
$ws.Range("J70").Value = 5097.9
$ws.Range("K70").Value = 6137.25
$ws.Range("L70").Value = 5097.9
$ws.Range("M70").Value = -5867.25
$ws.Range("N70").Value = -5637.9
$ws.Range("H73").Value = 5889.7856
$ws.Range("I73").Value = 6137.25
$ws.Range("J73").Value = 5097.9
$ws.Range("K73").Value = 6137.25
$ws.Range("L73").Value = 5097.9
$ws.Range("M73").Value = -5201.25
$ws.Range("N73").Value = -6969.9
$ws.Range("H98").Value = 23547.666
$ws.Range("J98").Value = 23547.666
$ws.Range("L98").Value = 23547.666
$ws.Range("N98").Value = -29537.666
$ws.Range("H132").Value = 5890.737
$ws.Range("I132").Value = 6492
$ws.Range("K132").Value = 19476
$ws.Range("M132").Value = -16946

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4493.467
$ws.Range("I40").Value = 3290.2
$ws.Range("J40").Value = 6900
$ws.Range("K40").Value = 3290.2
$ws.Range("L40").Value = 6900
$ws.Range("M40").Value = -3154.2
$ws.Range("N40").Value = -7172
$ws.Range("H82").Value = 4571.4287
$ws.Range("I82").Value = 4000
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 4000
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -3639
$ws.Range("N82").Value = -5722
$ws.Range("H85").Value = 4571.4287
$ws.Range("I85").Value = 4000
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 4000
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -2752
$ws.Range("N85").Value = -7496
$ws.Range("H132").Value = 17649.588
$ws.Range("I132").Value = 22920.334
$ws.Range("J132").Value = 4999.8
$ws.Range("K132").Value = 68761.00199999999
$ws.Range("L132").Value = 14999.4
$ws.Range("M132").Value = -66231.00199999999
$ws.Range("N132").Value = -20059.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 70015
$ws.Range("J22").Value = 70015
$ws.Range("L22").Value = 70015
$ws.Range("N22").Value = -70601
$ws.Range("H113").Value = 500.21738
$ws.Range("I113").Value = 327
$ws.Range("J113").Value = 825
$ws.Range("K113").Value = 981
$ws.Range("L113").Value = 2475
$ws.Range("M113").Value = 1189
$ws.Range("N113").Value = -6815
$ws.Range("H126").Value = 1025.5
$ws.Range("I126").Value = 796.0909
$ws.Range("J126").Value = 1866.6666
$ws.Range("K126").Value = 2388.2727
$ws.Range("L126").Value = 5599.9998
$ws.Range("M126").Value = 81.72730000000001
$ws.Range("N126").Value = -10539.9998
$ws.Range("H132").Value = 2552.7
$ws.Range("I132").Value = 1626.5555
$ws.Range("J132").Value = 3941.9167
$ws.Range("K132").Value = 4879.666499999999
$ws.Range("L132").Value = 11825.7501
$ws.Range("M132").Value = -2349.666499999999
$ws.Range("N132").Value = -16885.7501
$ws.Range("H136").Value = 2287.1936
$ws.Range("I136").Value = 1638.0476
$ws.Range("J136").Value = 3650.4
$ws.Range("K136").Value = 4914.142800000001
$ws.Range("L136").Value = 10951.2
$ws.Range("M136").Value = -2364.142800000001
$ws.Range("N136").Value = -16051.2
